$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '68.118.62'
$ws.Range("E2").Value = '  +0.00%  '
$ws.Range("D3").Value = '3.615.29'
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.30%  '
$ws.Range("D5").Value = "'586.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.20%  '
$ws.Range("D6").Value = "'193.44"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.69%  '
$ws.Range("D7").Value = '3.610.06'
$ws.Range("E7").Value = '  -1.26%  '
$ws.Range("D8").Value = "'0.620"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("E9").Value = '  +0.03%  '
$ws.Range("D10").Value = "'0.679"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -2.38%  '
$ws.Range("E11").Value = '  -0.10%  '
$ws.Range("D12").Value = "'55.33"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -1.79%  '
$ws.Range("D13").Value = "'0.0000291"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.86%  '
$ws.Range("D14").Value = "'10.02"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.70%  '
$ws.Range("D15").Value = '4.181.56'
$ws.Range("E15").Value = '  -1.54%  '
$ws.Range("D16").Value = '3.616.60'
$ws.Range("E16").Value = '  -1.28%  '
$ws.Range("E17").Value = '  -0.27%  '
$ws.Range("D18").Value = "'12.53"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.25%  '
$ws.Range("D19").Value = '67.947.32'
$ws.Range("E19").Value = '  +0.03%  '
$ws.Range("D20").Value = "'18.55"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -1.47%  '
$ws.Range("E21").Value = '  -1.94%  '
$ws.Range("D22").Value = "'405.22"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.58%  '
$ws.Range("D23").Value = "'13.26"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +22.29%  '
$ws.Range("D24").Value = "'4.26"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.71%  '
$ws.Range("D25").Value = "'86.08"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -1.90%  '
$ws.Range("E26").Value = '  +0.72%  '
$ws.Range("B27").Value = 'Toncoin'
$ws.Range("C27").Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range("D27").Value = "'3.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +6.81%  '
$ws.Range("B28").Value = 'InternetComputer(DFINITY)'
$ws.Range("C28").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D28").Value = "'12.60"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +1.19%  '
$ws.Range("E29").Value = '  +0.69%  '
$ws.Range("D30").Value = "'8.10"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +13.44%  '
$ws.Range("D31").Value = "'9.17"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.13%  '
$ws.Range("D32").Value = "'31.64"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -0.43%  '
$ws.Range("D33").Value = "'686.45"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +13.25%  '
$ws.Range("D34").Value = "'12.26"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.49%  '
$ws.Range("E35").Value = '  +2.30%  '
$ws.Range("D36").Value = "'64.59"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -5.14%  '
$ws.Range("D37").Value = "'42.78"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -3.04%  '
$ws.Range("D38").Value = "'0.423"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.61%  '
$ws.Range("D39").Value = "'0.999"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.09%  '
$ws.Range("D40").Value = '0.0₃0790'
$ws.Range("E40").Value = '  +4.08%  '
$ws.Range("D41").Value = "'2.95"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +18.05%  '
$ws.Range("B42").Value = 'ThetaToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta'
$ws.Range("D42").Value = "'3.13"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +9.12%  '
$ws.Range("B43").Value = 'Maker'
$ws.Range("C43").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D43").Value = '3.186.79'
$ws.Range("E43").Value = '  +15.57%  '
$ws.Range("E44").Value = '  -0.30%  '
$ws.Range("E45").Value = '  -0.30%  '
$ws.Range("D46").Value = "'0.0423"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.02%  '
$ws.Range("E47").Value = '  -2.32%  '
$ws.Range("B48").Value = 'THORChain'
$ws.Range("C48").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D48").Value = "'8.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.46%  '
$ws.Range("B49").Value = 'ApeXProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D49").Value = "'3.10"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -4.36%  '
$ws.Range("D50").Value = "'142.96"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.28%  '
$ws.Range("E51").Value = '  -0.59%  '
